$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "27.101.98"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "'" + "1.892.08"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("D4").Value = "'" + "1.000"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'" + "306.73"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "'" + "1.000"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'" + "0.5190"
$ws.Range("E7").Value = "  +2.84%  "

$ws.Range("D8").Value = "'" + "0.3761"
$ws.Range("E8").Value = "  +3.22%  "

$ws.Range("D9").Value = "'" + "0.07222"
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("D10").Value = "'" + "21.16"
$ws.Range("E10").Value = "  +2.67%  "

$ws.Range("D11").Value = "'" + "0.9029"
$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'" + "0.07661"
$ws.Range("E12").Value = "  +1.96%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'" + "1.915.81"
$ws.Range("E13").Value = "  +2.84%  "

$ws.Range("D14").Value = "'" + "94.47"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").Value = "'" + "5.242"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "'" + "1.001"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'" + "0.000008512"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("E18").Value = "  +1.67%  "

$ws.Range("D19").Value = "'" + "0.9996"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'" + "27.146.57"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").Value = "'" + "5.066"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").Value = "'" + "2.182.73"
$ws.Range("E22").Value = "  +4.14%  "

$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").Value = "'" + "6.387"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "'" + "2.309"
$ws.Range("E25").Value = "  +11.59%  "

$ws.Range("D26").Value = "'" + "145.68"
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").Value = "'" + "18.07"
$ws.Range("E27").Value = "  +1.07%  "

$ws.Range("D28").Value = "'" + "1.728"
$ws.Range("E28").Value = "  -2.98%  "

$ws.Range("D29").Value = "'" + "114.49"
$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").Value = "'" + "4.958"
$ws.Range("E30").Value = "  +6.47%  "

$ws.Range("D31").Value = "'" + "4.800"
$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("D32").Value = "'" + "0.09212"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").Value = "'" + "0.05059"
$ws.Range("E33").Value = "  -1.52%  "

$ws.Range("D34").Value = "'" + "1.248"
$ws.Range("E34").Value = "  +8.31%  "

$ws.Range("D35").Value = "'" + "0.7740"
$ws.Range("E35").Value = "  +3.45%  "

$ws.Range("D36").Value = "'" + "2.987"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Value = "'" + "3.283"
$ws.Range("E37").Value = "  +2.79%  "

$ws.Range("D38").Value = "'" + "2.603"
$ws.Range("E38").Value = "  +1.68%  "

$ws.Range("D39").Value = "'" + "0.5677"
$ws.Range("E39").Value = "  +1.89%  "

$ws.Range("D40").Value = "'" + "0.01993"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").Value = "'" + "1.073"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").Value = "'" + "9.053"
$ws.Range("E42").Value = "  +6.11%  "

$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("D44").Value = "'" + "119.31"
$ws.Range("E44").Value = "  +3.03%  "

$ws.Range("D45").Value = "'" + "0.1512"
$ws.Range("E45").Value = "  +2.95%  "

$ws.Range("D46").Value = "'" + "0.4861"
$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "'" + "10.14"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("D49").Value = "'" + "1.598"
$ws.Range("E49").Value = "  +2.97%  "

$ws.Range("D50").Value = "'" + "37.70"
$ws.Range("E50").Value = "  +2.67%  "

$ws.Range("D51").Value = "'" + "63.99"
$ws.Range("E51").Value = "  +1.62%  "
